$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 10) with the latest profit allocation data
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "09/11/2025"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = 0.1236841051452084
$ws.Range("C10").Value = 0.8763158948547916
